$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 12:50:41"
$wsZhCn.Range("H2").Value = "2016-03-19 12:50:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 12:50:44"
$wsDeDe.Range("H2").Value = "2016-03-19 12:51:05"
